# Update cryptocurrency price/volume data in the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.640.06"
$ws.Range("D3").Value = "2.437.63"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  +0.10%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "567.14"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "145.50"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.69%  "
$ws.Range("E7").Value = "  -0.10%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.532"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  +2.05%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.355"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("E13").Value = "  +6.21%  "
$ws.Range("E14").Value = "  +5.83%  "
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").Value = "62.562.87"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "2.448.20"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("E19").Value = "  +2.82%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "323.96"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("E22").Value = "  -0.06%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.84"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +8.20%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "67.33"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("E25").Value = "  -0.78%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "582.08"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.91%  "
$ws.Range("E27").Value = "  +9.42%  "
$ws.Range("E29").Value = "  +0.15%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.43"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +3.96%  "
$ws.Range("E31").Value = "  +4.99%  "
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  +1.30%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "18.77"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("E39").Value = "  -0.07%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "148.13"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("E41").Value = "  +2.78%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +10.78%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "148.27"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("E47").Value = "  +4.33%  "
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("E49").Value = "  +3.55%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0921"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.66%  "
